$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (B1:G1) from abbreviated lowercase labels to
# capitalised "Stats" section headers.
$ws.Range("B1").Value = "Pld"
$ws.Range("C1").Value = "Won"
$ws.Range("D1").Value = "Lost"
$ws.Range("E1").Value = "Tied"
$ws.Range("F1").Value = "Net RR"
$ws.Range("G1").Value = "Pts"

# Move the selection to G1 to match the saved view state.
$ws.Range("G1").Select()
